$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Foundation-1"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Foundation-2"

$ws1.Range("A1:C18").Copy($ws2.Range("A1:C18"))

$ws2.Activate() | Out-Null
$ws2.Range("E4").Select() | Out-Null
